$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data rows (rows 2-10), keep header row 1 (and its styles) intact
$ws.Range("A2:T10").Clear()

# Write the text columns (A-D) column-by-column so the shared-string table is
# rebuilt in the same first-seen order as the refreshed NATMI export

# Column A
$ws.Range("A2").Value = "FAPs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "MuSCs"
$ws.Range("A6").Value = "MuSCs"
$ws.Range("A7").Value = "MuSCs"

# Column B
$ws.Range("B2").Value = "Fgf18"
$ws.Range("B3").Value = "Fgf18"
$ws.Range("B4").Value = "Fgf18"
$ws.Range("B5").Value = "Fgf18"
$ws.Range("B6").Value = "Fgf18"
$ws.Range("B7").Value = "Fgf18"

# Column C
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("C7").Value = "Fgfr2"

# Column D
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"

# Write the numeric columns (E-T) row-by-row

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.71557066666667
$ws.Range("H2").Value = 32.146712
$ws.Range("I2").Value = 0.9375025736567436
$ws.Range("J2").Value = 0.9375025736567436
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05601
$ws.Range("N2").Value = 0.16803
$ws.Range("O2").Value = 0.02710547761971223
$ws.Range("P2").Value = 0.02710547761971223
$ws.Range("Q2").Value = 0.60017911304
$ws.Range("R2").Value = 5.401612017360001
$ws.Range("S2").Value = 0.02541145502867548
$ws.Range("T2").Value = 0.02541145502867548

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.71557066666667
$ws.Range("H3").Value = 32.146712
$ws.Range("I3").Value = 0.9375025736567436
$ws.Range("J3").Value = 0.9375025736567436
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.864751
$ws.Range("N3").Value = 5.594253
$ws.Range("O3").Value = 0.902427539668559
$ws.Range("P3").Value = 0.9024275396685592
$ws.Range("Q3").Value = 19.98187111623733
$ws.Range("R3").Value = 179.836840046136
$ws.Range("S3").Value = 0.8460281409779972
$ws.Range("T3").Value = 0.8460281409779973

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.71557066666667
$ws.Range("H4").Value = 32.146712
$ws.Range("I4").Value = 0.9375025736567436
$ws.Range("J4").Value = 0.9375025736567436
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.145611
$ws.Range("N4").Value = 0.436833
$ws.Range("O4").Value = 0.07046698271172858
$ws.Range("P4").Value = 0.07046698271172858
$ws.Range("Q4").Value = 1.560304960344
$ws.Range("R4").Value = 14.042744643096
$ws.Range("S4").Value = 0.0660629776500708
$ws.Range("T4").Value = 0.0660629776500708

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.71434
$ws.Range("H5").Value = 2.14302
$ws.Range("I5").Value = 0.06249742634325634
$ws.Range("J5").Value = 0.06249742634325634
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05601
$ws.Range("N5").Value = 0.16803
$ws.Range("O5").Value = 0.02710547761971223
$ws.Range("P5").Value = 0.02710547761971223
$ws.Range("Q5").Value = 0.0400101834
$ws.Range("R5").Value = 0.3600916506
$ws.Range("S5").Value = 0.001694022591036748
$ws.Range("T5").Value = 0.001694022591036749

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.71434
$ws.Range("H6").Value = 2.14302
$ws.Range("I6").Value = 0.06249742634325634
$ws.Range("J6").Value = 0.06249742634325634
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.864751
$ws.Range("N6").Value = 5.594253
$ws.Range("O6").Value = 0.902427539668559
$ws.Range("P6").Value = 0.9024275396685592
$ws.Range("Q6").Value = 1.33206622934
$ws.Range("R6").Value = 11.98859606406
$ws.Range("S6").Value = 0.05639939869056181
$ws.Range("T6").Value = 0.05639939869056182

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.71434
$ws.Range("H7").Value = 2.14302
$ws.Range("I7").Value = 0.06249742634325634
$ws.Range("J7").Value = 0.06249742634325634
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.145611
$ws.Range("N7").Value = 0.436833
$ws.Range("O7").Value = 0.07046698271172858
$ws.Range("P7").Value = 0.07046698271172858
$ws.Range("Q7").Value = 0.10401576174
$ws.Range("R7").Value = 0.9361418556600001
$ws.Range("S7").Value = 0.004404005061657775
$ws.Range("T7").Value = 0.004404005061657775
